$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Friday 09/05/2025 - "Trying to figure out shooting" ---
$ws.Range("A12").Value = "Friday"
$ws.Range("B12").Value = 45786
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Value = 0.70833333333333337
$ws.Range("F12").Value = "Trying to figure out shooting"

# --- Row 13: Sunday 11/05/2025 - "Working on smoke dash" ---
$ws.Range("A13").Value = "Sunday"
$ws.Range("B13").Value = 45788
$ws.Range("C13").Value = 0.625
$ws.Range("D13").Value = 0.70833333333333337
$ws.Range("F13").Value = "Working on smoke dash"

# --- Extend the shared "time spent" formula (D-C) down through the new rows ---
$ws.Range("E6:E13").Formula = "=D6-C6"

# --- Copy the formatting of an existing data row onto the two new rows ---
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F12").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:F13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Move the active selection to the new last row, like the author's save ---
$ws.Range("F13").Select() | Out-Null

# --- Recalculate so dependent totals (I1) pick up the new rows ---
$excel.Calculate()
